# Non-Oncology Excel reports comparison
# Inserts a "StandardExcelReport-...-2023_" row ahead of the existing
# ExcelReport-/WordReport- rows for each of the four scenario blocks
# (Clinical, Economic, Quality of Life, Real-world Evidence) in column K,
# shifting the remaining report names down by one row per block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clinical block (rows 2-4) ---
$ws.Range("K2").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-2023_"
$ws.Range("K3").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Clinical-"
$ws.Range("K4").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"

# --- Economic block (rows 5-7) ---
$ws.Range("K5").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Economic-2023_"
$ws.Range("K6").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Economic-"
$ws.Range("K7").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Economic-"

# --- Quality of Life block (rows 8-10) ---
$ws.Range("K8").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-2023_"
$ws.Range("K9").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Quality of Life-"
$ws.Range("K10").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"

# --- Real-world Evidence block (rows 11-13) ---
$ws.Range("K11").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-2023_"
$ws.Range("K12").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Real-world Evidence-"
$ws.Range("K13").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"

# Update selection / view to match the reviewed range
$ws.Range("K2:K13").Select()
